# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
#
# The underlying data rows for several fixtures were re-sequenced: the
# match record that used to sit on one row now sits on its paired row
# (and vice versa). Column A (the running `id` sequence) stays put since
# it always mirrors the row's position; columns B..AC (id/B, HomeTeam,
# AwayTeam, scores, odds, P/L, etc.) swap wholesale between each pair of
# rows. C/D/E (Div, Div Original Name, Date) are identical between the
# two rows of a pair, so including them in the swapped range is harmless.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$RowA,
        [int]$RowB
    )

    $rangeA = $ws.Range("B$RowA`:AC$RowA")
    $rangeB = $ws.Range("B$RowB`:AC$RowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

Swap-RowData 213 214
Swap-RowData 215 216
Swap-RowData 229 231
Swap-RowData 232 233
Swap-RowData 251 252
Swap-RowData 263 265
